$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

function CopyFormat($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats) | Out-Null
    $excel.CutCopyMode = $false
}

# --- Step 1: retarget the existing "Date of clinic visit" note in row 16 to a
# new, more specific note about linking the competency to the three clinical
# cases. This is done in place, before any rows move. ---
$ws.Range("C16").Value = "Link this competency to the three clinical cases"

# Row 16 picks up the lighter-weight grid formatting used by the row above it
# (no longer the "last sub-row" of the competency-6 block).
CopyFormat "D3" "D14"
CopyFormat "D4" "D15"
CopyFormat "A2" "A16"
CopyFormat "B2" "B16"
CopyFormat "B2" "C16"
CopyFormat "D4" "D16"

# --- Step 2: insert a new row to hold the original "Date of clinic visit"
# note (now pushed one row down). ---
$ws.Rows(17).Insert()

$ws.Range("C17").Value = "Date of clinic visit"

# Give the freshly-inserted row the "last sub-row" grid formatting that row 16
# used to have.
CopyFormat "A10" "A17"
CopyFormat "B8" "B17"
CopyFormat "B2" "C17"
CopyFormat "D4" "D17"

# Restore the author's final cursor position.
$ws.Range("F17").Select() | Out-Null
